# Changing names for delivery forms
# Replace the "skate size" question with a "hygiene kit supply" yes/no question,
# and normalize the yes_no choice list's data_value casing (true/false -> TRUE/FALSE).

$wb = $excel.ActiveWorkbook

# --- survey sheet: the skate-size instruction/question becomes a hygiene-kit question ---
$survey = $wb.Worksheets.Item("survey")

# Row 3: instruction text shown to the user
$survey.Range("F3").Value = "Distribute the hygiene kit"

# Row 4: was a blank-type/blank-values_list "decimal" skate-size question,
# now a select_one yes_no question asking about supplies lasting the month.
$survey.Range("F4").Value = "Did the supplies last you the entire month?"
$survey.Range("E4").Value = "supply"
$survey.Range("C4").Value = "select_one"
$survey.Range("D4").Value = "yes_no"
$survey.Range("G4").Value = ""

# --- model sheet: the calculated field's declared type/name changes to match ---
$model = $wb.Worksheets.Item("model")
$model.Range("A5").Value = "string"
$model.Range("B5").Value = "supply"

# --- choices sheet: yes_no data values are now uppercase TRUE/FALSE ---
# (leading apostrophe forces these to stay text cells instead of Excel
# auto-converting the literal words TRUE/FALSE into boolean cells)
$choices = $wb.Worksheets.Item("choices")
$choices.Range("B2").Value = "'TRUE"
$choices.Range("B3").Value = "'FALSE"
